$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "49+45=94",
    "5+58=63",
    "5+56=61",
    "58+13=71",
    "60-22=38",
    "28+18=46",
    "60-49=11",
    "14+78=92",
    "77+9=86",
    "55+8=63",
    "7+38=45",
    "53+38=91",
    "80-6=74",
    "9+15=24",
    "29+43=72",
    "36-29=7",
    "47+45=92",
    "80-62=18",
    "46+26=72",
    "16+69=85",
    "96-29=67",
    "37+35=72",
    "65-39=26",
    "26-19=7",
    "81-72=9",
    "4+39=43",
    "91-44=47",
    "16+37=53",
    "50-19=31",
    "37-29=8",
    "6+18=24",
    "49+15=64",
    "90-65=25",
    "70-58=12",
    "27+7=34",
    "84-48=36",
    "56-47=9",
    "23-6=17",
    "27+7=34",
    "6+7=13",
    "39+17=56",
    "19+52=71",
    "52-17=35",
    "92-38=54",
    "96-88=8",
    "8+75=83",
    "18+9=27",
    "15+76=91",
    "14+8=22",
    "86-29=57",
    "58+29=87",
    "36+19=55",
    "8+53=61",
    "52-17=35",
    "23+19=42",
    "64+28=92",
    "6+88=94",
    "28+57=85",
    "97-58=39",
    "52-26=26",
    "68+18=86",
    "15+9=24",
    "4+59=63",
    "38+27=65",
    "25+6=31",
    "30-4=26",
    "71-58=13",
    "66-8=58",
    "25+19=44",
    "66-57=9",
    "92-79=13",
    "49+26=75",
    "84-66=18",
    "64-55=9",
    "92-28=64",
    "13-7=6",
    "25+36=61",
    "18+75=93",
    "12+39=51",
    "79+12=91",
    "63-56=7",
    "93-59=34",
    "80-71=9",
    "9+74=83",
    "67+24=91",
    "98-79=19",
    "33-7=26",
    "11-7=4",
    "53-14=39",
    "80-44=36",
    "81-24=57",
    "90-52=38",
    "63-46=17",
    "86-59=27",
    "17+7=24",
    "27+54=81",
    "5+29=34",
    "94-45=49",
    "39+7=46",
    "83-49=34"
)

$rows = 20
$cols = 5
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i = $i + 1
    }
}
Write-Output "done: $i cells updated"